$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert the two new paragraphs in front of the "Please note: ..."
#    paragraph:
#      a) a bold+underlined paragraph (wrapped in a bookmark) with the
#         new "trial bundle" warning text
#      b) a blank spacer paragraph
#    We anchor on the end of the paragraph that currently precedes the
#    "Please note" paragraph (the empty "<additionalInfo>" spacer) and
#    insert raw WordprocessingML there, which is the most reliable way
#    to reproduce the exact paragraph/run formatting from the diff.
# ------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Please note:*") {
        $target = $para
        break
    }
}

$anchor = $d.Range($target.Range.Start, $target.Range.Start)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="5" w:name="_Hlk150165830"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>The Claimant (or in the event they act in person and the Defendant is represented, the Defendant) must bring to court for the start of the trial a paper copy of the electronic trial bundle for use by witnesses. A failure to do so may result in the imposition of sanctions.</w:t>
  </w:r>
  <w:bookmarkEnd w:id="5"/>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$anchor.InsertXML($xml)

# ------------------------------------------------------------------
# 2) Merge the two runs of the "Please note" paragraph's second
#    sentence into a single run and drop the lastRenderedPageBreak
#    that used to sit between them — simplest reliable way is a
#    Find/Replace across the exact (unchanged) sentence, which
#    Word always writes back out as one run.
# ------------------------------------------------------------------

$sentence = "Cases are listed in accordance with local hearing arrangements determined by the Judiciary and implemented by the court staff. Every effort is made to ensure that hearings start at the time specified. However, listing practices or other factors may mean that you experience a delay, an adjournment at short notice or your case may be released to a different court hearing centre, in which case you will be notified."

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)
